$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "procedural creation of platforms" feature row entirely.
# This shifts every subsequent row up by one and auto-adjusts the
# table range, autofilter, sort state and data validation sqrefs.
$ws.Rows(15).Delete()

# Update two of the remaining feature descriptions.
$ws.Range("A16").Value = "Death zones that are attached to objects"
$ws.Range("A17").Value = "Implemetation of obstacles and dead zones to cross"

# Fill in the "Assigned To" column with the team member responsible
# for each feature.
$ws.Range("C8").Value  = "Jacob"
$ws.Range("C9").Value  = "Brennan"
$ws.Range("C10").Value = "Brennan"
$ws.Range("C11").Value = "Jacob"
$ws.Range("C12").Value = "Jacob"
$ws.Range("C13").Value = "Brennan"
$ws.Range("C14").Value = "Brennan"
$ws.Range("C15").Value = "Jacob"
$ws.Range("C16").Value = "Jacob"
$ws.Range("C17").Value = "Brennan"
$ws.Range("C18").Value = "Both?"
$ws.Range("C19").Value = "Jacob"
$ws.Range("C20").Value = "Brennan"
$ws.Range("C21").Value = "Brennan"
$ws.Range("C22").Value = "Jacob"
$ws.Range("C23").Value = "Both?"

# Restore the view's selection to match the saved state.
$ws.Range("C21").Select()
